$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2 Backlog")

# Row 9 ("Refactoring Code") gains actual-effort / status updates
$ws.Range("I9").Value = "6h"
$ws.Range("J9").Value = "8h"
$ws.Range("K9").Value = "erledigt"

# Drop the two trailing blank formatting-only rows (11 and 12)
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(11).Delete()

# Move the active selection to reflect the edited area
$ws.Range("I10").Select()
